# On slide 23 (the "Tokens" slide, sldId 279), the "Identifier:" content
# placeholder currently holds the explanatory sentence as a single plain
# run:
#   "An identifier is a name given to a variable, function, or other
#    entity in a program. Will follow rules, ..."
#
# The edit bolds a few key phrases within that sentence ("name given to a
# variable", "function", "other entity"). In OOXML this is realized by
# splitting the run into several runs that keep the same rPr except for
# b="1" on the bolded spans, which is exactly what happens when you set
# Font.Bold on a sub-range of an existing run via the PowerPoint object
# model.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(23)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

$tr.Find("name given to a variable", 0).Font.Bold = $true
$tr.Find("function", 0).Font.Bold = $true
$tr.Find("other entity ", 0).Font.Bold = $true
